$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the real id value in A3 with a string that doesn't correspond to
# any existing record, turning this row into the "import should fail" case.
$ws.Range("A3").Value = "ThatIdDoesntExist"

# Update the selection / active cell shown when the sheet is reopened.
$ws.Range("A2").Select()
